$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Published Values")

# --- Append the new row (row 14) of published-rate data ---
# All source cells in this sheet are stored as text, so force the row to
# Text format before writing the values (keeps "697.85" etc. as a string
# instead of being auto-coerced to a number), then clear the format again
# so no stray style index gets attached to the new cells.
$rowIndex = 14
$values = @(
    "2026-01-02",
    "2026-01-02 21:19:48",
    "697.85",
    "697.85",
    "700.79",
    "700.79",
    "702.88",
    "2026/01/02 21:19:48",
    "2026-01-02 13:41:27",
    "https://www.bankofchina.com/sourcedb/whpj/enindex_1619.html"
)

$newRowRange = $ws.Range("A14:J14")
$newRowRange.NumberFormat = "@"
for ($col = 1; $col -le $values.Length; $col++) {
    $ws.Cells.Item($rowIndex, $col).Value = $values[$col - 1]
}
$newRowRange.ClearFormats()

# --- Grow the AutoFilter range to cover the new row ---
$ws.AutoFilterMode = $false
$ws.Range("A1:J14").AutoFilter()

# --- Keep the hidden _FilterDatabase defined name in sync ---
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -eq "All Published Values!_FilterDatabase") {
        $n.RefersTo = "='All Published Values'!`$A`$1:`$J`$14"
    }
}

# --- Update the Daily Summary publishes count for 2026-01-02 ---
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Range("B4").Value = 13
